# Usman - New users added to cobaltUsers.xls
# Adds 33 new user rows (81-113) to the "Users" sheet, matching rows already
# present in columns A (username), B (password), E/F (notes), G (email with
# mailto hyperlink), then updates the sheet view / selection / dimension
# accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Data for the new rows: row number, username (col A), email (col G).
# Column B is always "Password1", column E is always the standing warning
# text, column F is always "N" - these already exist as shared strings so no
# new strings get created for them.
# ---------------------------------------------------------------------------
$rowsData = @(
    @{Row=81;  User='KHPaddUser1';            Email='KHPaddUser1@mailinator.com'}
    @{Row=82;  User='KHPaddUser2';            Email='KHPaddUser2@mailinator.com'}
    @{Row=83;  User='KHPaddUser3';            Email='KHPaddUser3@mailinator.com'}
    @{Row=84;  User='KHPaddUser4';            Email='KHPaddUser4@mailinator.com'}
    @{Row=85;  User='KHPaddUser5';            Email='KHPaddUser5@mailinator.com'}
    @{Row=86;  User='KHPaddUser6';            Email='KHPaddUser6@mailinator.com'}
    @{Row=87;  User='SearchWhatsMarketUser1'; Email='SearchWhatsMarketUser1@mailinator.com '}
    @{Row=88;  User='SearchWhatsMarketUser2'; Email='SearchWhatsMarketUser2@mailinator.com '}
    @{Row=89;  User='SearchWhatsMarketUser3'; Email='SearchWhatsMarketUser3@mailinator.com '}
    @{Row=90;  User='SearchWhatsMarketUser4'; Email='SearchWhatsMarketUser4@mailinator.com '}
    @{Row=91;  User='SearchWhatsMarketUser5'; Email='SearchWhatsMarketUser5@mailinator.com '}
    @{Row=92;  User='SearchWhatsMarketUser6'; Email='SearchWhatsMarketUser6@mailinator.com '}
    @{Row=93;  User='SearchWhatsMarketUser7'; Email='SearchWhatsMarketUser7@mailinator.com '}
    @{Row=94;  User='SearchWhatsMarketUser8'; Email='SearchWhatsMarketUser8@mailinator.com '}
    @{Row=95;  User='SearchKnowHowUser1';     Email='SearchKnowHowUser1@mailinator.com '}
    @{Row=96;  User='SearchKnowHowUser2';     Email='SearchKnowHowUser2@mailinator.com '}
    @{Row=97;  User='SearchKnowHowUser3';     Email='SearchKnowHowUser3@mailinator.com '}
    @{Row=98;  User='SearchKnowHowUser4';     Email='SearchKnowHowUser4@mailinator.com '}
    @{Row=99;  User='SearchKnowHowUser5';     Email='SearchKnowHowUser5@mailinator.com '}
    @{Row=100; User='SearchKnowHowUser6';     Email='SearchKnowHowUser6@mailinator.com '}
    @{Row=101; User='SearchKnowHowUser7';     Email='SearchKnowHowUser7@mailinator.com '}
    @{Row=102; User='SearchKnowHowUser8';     Email='SearchKnowHowUser8@mailinator.com '}
    @{Row=103; User='AskUser1';               Email='AskUser1@mailinator.com '}
    @{Row=104; User='AskUser2';               Email='AskUser2@mailinator.com '}
    @{Row=105; User='AskUser3';               Email='AskUser3@mailinator.com '}
    @{Row=106; User='AskUser4';               Email='AskUser4@mailinator.com '}
    @{Row=107; User='AskUser5';               Email='AskUser5@mailinator.com '}
    @{Row=108; User='AskUser6';               Email='AskUser6@mailinator.com '}
    @{Row=109; User='AssetPageUser1';         Email='AssetPageUser1@mailinator.com '}
    @{Row=110; User='AssetPageUser2';         Email='AssetPageUser2@mailinator.com '}
    @{Row=111; User='AssetPageUser3';         Email='AssetPageUser3@mailinator.com '}
    @{Row=112; User='AssetPageUser4';         Email='AssetPageUser4@mailinator.com '}
    @{Row=113; User='AssetPageUser5';         Email='AssetPageUser5@mailinator.com '}
)

$passwordText = 'Password1'
$noteText     = 'THIS IS IN USE 24/7 - DO NOT USE!'
$flagText     = 'N'

# ---------------------------------------------------------------------------
# Write the cell values first (rows 81-113, columns A/B/E/F/G). D/C are left
# blank, exactly like the pre-existing rows just above (74-80).
# ---------------------------------------------------------------------------
foreach ($r in $rowsData) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.User
    $ws.Range("B$row").Value = $passwordText
    $ws.Range("E$row").Value = $noteText
    $ws.Range("F$row").Value = $flagText
    $ws.Range("G$row").Value = $r.Email
}

# ---------------------------------------------------------------------------
# Formatting:
#  - Column A on the new rows uses a 10pt Arial font (a style that didn't
#    exist before this edit), built once on a scratch cell then stamped onto
#    every new A-cell via a format-only paste so only a single new font/style
#    pair gets added to the workbook instead of one per cell.
#  - Columns E/F/G on the new rows copy the formatting already used by the
#    identical columns in row 80 (the previous last row), reusing the
#    existing styles rather than creating new ones.
# ---------------------------------------------------------------------------
$template = $ws.Range("Z1000")
$template.Font.Name = "Arial"
$template.Font.Size = 10
$template.Copy()
$ws.Range("A81:A113").PasteSpecial(-4122)
$template.Clear()

$ws.Range("E80:G80").Copy()
$ws.Range("E81:G113").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Hyperlinks on column G. Every new row except G81 and G112 gets a mailto:
# hyperlink to its own address (added in this exact order so the relationship
# ids line up: G82..G110, then G113, then G111 last). G111's hyperlink keeps
# a stale display caption copied from AssetPageUser1 (row 109) even though
# the cell text/address is AssetPageUser3.
# ---------------------------------------------------------------------------
$hyperlinkOrder = @(82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102,103,104,105,106,107,108,109,110,113,111)
$displayOverride = @{ 111 = 'AssetPageUser1@mailinator.com ' }

$byRow = @{}
foreach ($r in $rowsData) { $byRow[$r.Row] = $r }

foreach ($row in $hyperlinkOrder) {
    $info = $byRow[$row]
    $address = "mailto:" + $info.Email.Trim()
    if ($displayOverride.ContainsKey($row)) {
        $ws.Hyperlinks.Add($ws.Range("G$row"), $address, "", "", $displayOverride[$row]) | Out-Null
    } else {
        $ws.Hyperlinks.Add($ws.Range("G$row"), $address) | Out-Null
    }
}

# ---------------------------------------------------------------------------
# Sheet view: scrolled down and selection now covers the newly added block.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 67
$ws.Range("A81:G113").Select()
